$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.092262
$ws.Range("H2").Value = 0.276786
$ws.Range("I2").Value = 0.674176843971804
$ws.Range("J2").Value = 0.6741768439718039
$ws.Range("M2").Value = 9.841031333333333
$ws.Range("N2").Value = 29.523094
$ws.Range("O2").Value = 0.1083017349730097
$ws.Range("P2").Value = 0.1125970533891552
$ws.Range("Q2").Value = 0.9079532328759999
$ws.Range("R2").Value = 8.171579095883999
$ws.Range("S2").Value = 0.07301452188077442
$ws.Range("T2").Value = 0.07591032609442534

$ws.Range("G3").Value = 0.092262
$ws.Range("H3").Value = 0.276786
$ws.Range("I3").Value = 0.674176843971804
$ws.Range("J3").Value = 0.6741768439718039
$ws.Range("O3").Value = 0.3504595127507141
$ws.Range("P3").Value = 0.3643589687437936
$ws.Range("Q3").Value = 2.938095568584
$ws.Range("R3").Value = 26.442860117256
$ws.Range("S3").Value = 0.2362716882461726
$ws.Range("T3").Value = 0.2456423796205119

$ws.Range("G4").Value = 0.092262
$ws.Range("H4").Value = 0.276786
$ws.Range("I4").Value = 0.674176843971804
$ws.Range("J4").Value = 0.6741768439718039
$ws.Range("M4").Value = 16.16670066666667
$ws.Range("N4").Value = 48.500102
$ws.Range("O4").Value = 0.1779164877830196
$ws.Range("P4").Value = 0.1849727733236046
$ws.Range("Q4").Value = 1.491572136908
$ws.Range("R4").Value = 13.424149232172
$ws.Range("S4").Value = 0.1199471762241041
$ws.Range("T4").Value = 0.1247043605400196

$ws.Range("G5").Value = 0.092262
$ws.Range("H5").Value = 0.276786
$ws.Range("I5").Value = 0.674176843971804
$ws.Range("J5").Value = 0.6741768439718039
$ws.Range("M5").Value = 10.399077
$ws.Range("N5").Value = 20.798154
$ws.Range("O5").Value = 0.1144430947397913
$ws.Range("P5").Value = 0.07932132236322763
$ws.Range("Q5").Value = 0.959439642174
$ws.Range("R5").Value = 5.756637853043999
$ws.Range("S5").Value = 0.07715488442603863
$ws.Range("T5").Value = 0.05347659877051087

$ws.Range("G6").Value = 0.092262
$ws.Range("H6").Value = 0.276786
$ws.Range("I6").Value = 0.674176843971804
$ws.Range("J6").Value = 0.6741768439718039
$ws.Range("M6").Value = 22.614852
$ws.Range("N6").Value = 67.84455600000001
$ws.Range("O6").Value = 0.2488791697534654
$ws.Range("P6").Value = 0.258749882180219
$ws.Range("Q6").Value = 2.086491475224
$ws.Range("R6").Value = 18.778423277016
$ws.Range("S6").Value = 0.1677885731947142
$ws.Range("T6").Value = 0.1744431789463361

$ws.Range("I7").Value = 0.3258231560281961
$ws.Range("J7").Value = 0.3258231560281961
$ws.Range("M7").Value = 9.841031333333333
$ws.Range("N7").Value = 29.523094
$ws.Range("O7").Value = 0.1083017349730097
$ws.Range("P7").Value = 0.1125970533891552
$ws.Range("Q7").Value = 0.4388050264657778
$ws.Range("R7").Value = 3.949245238192
$ws.Range("S7").Value = 0.03528721309223528
$ws.Range("T7").Value = 0.03668672729472983

$ws.Range("I8").Value = 0.3258231560281961
$ws.Range("J8").Value = 0.3258231560281961
$ws.Range("O8").Value = 0.3504595127507141
$ws.Range("P8").Value = 0.3643589687437936
$ws.Range("S8").Value = 0.1141878245045415
$ws.Range("T8").Value = 0.1187165891232817

$ws.Range("I9").Value = 0.3258231560281961
$ws.Range("J9").Value = 0.3258231560281961
$ws.Range("M9").Value = 16.16670066666667
$ws.Range("N9").Value = 48.500102
$ws.Range("O9").Value = 0.1779164877830196
$ws.Range("P9").Value = 0.1849727733236046
$ws.Range("Q9").Value = 0.7208624049262222
$ws.Range("R9").Value = 6.487761644336
$ws.Range("S9").Value = 0.05796931155891542
$ws.Range("T9").Value = 0.06026841278358497

$ws.Range("I10").Value = 0.3258231560281961
$ws.Range("J10").Value = 0.3258231560281961
$ws.Range("M10").Value = 10.399077
$ws.Range("N10").Value = 20.798154
$ws.Range("O10").Value = 0.1144430947397913
$ws.Range("P10").Value = 0.07932132236322763
$ws.Range("Q10").Value = 0.463687910712
$ws.Range("R10").Value = 2.782127464272
$ws.Range("S10").Value = 0.03728821031375263
$ws.Range("T10").Value = 0.02584472359271675

$ws.Range("I11").Value = 0.3258231560281961
$ws.Range("J11").Value = 0.3258231560281961
$ws.Range("M11").Value = 22.614852
$ws.Range("N11").Value = 67.84455600000001
$ws.Range("O11").Value = 0.2488791697534654
$ws.Range("P11").Value = 0.258749882180219
$ws.Range("Q11").Value = 1.008381174112
$ws.Range("R11").Value = 9.075430567008002
$ws.Range("S11").Value = 0.08109059655875125
$ws.Range("T11").Value = 0.08430670323388284
